$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: bold formatting (A1:J1) ---
$ws.Range("A1:J1").Font.Bold = $true

# --- Column A width: narrower than before ---
$ws.Columns.Item(1).ColumnWidth = 20.5

# --- Row 8: Alpha, Beta Attention / Unbalanced ---
$ws.Range("A8").Value = "Alpha, Beta Attention"
$ws.Range("A8").WrapText = $true
$ws.Range("B8").Value = "Alpha, Beta Attention"
$ws.Range("C8").Value = "Unbalanced"
$ws.Range("D8").Value = 431.98700000000002
$ws.Range("E8").Value = 0.001
$ws.Range("F8").Value = 10
$ws.Range("G8").Value = 0.67741935483870896
$ws.Range("H8").Value = 0.458515283842794
$ws.Range("I8").Value = 0.546875
$ws.Range("J8").Value = 0.882634375062886

# --- Row 9: Alpha, Beta Attention / Balanced (bold-highlighted best result) ---
$ws.Range("A9").Value = "Alpha, Beta Attention"
$ws.Range("A9").WrapText = $true
$ws.Range("B9").Value = "Alpha, Beta Attention"
$ws.Range("C9").Value = "Balanced"
$ws.Range("D9").Value = 852.77
$ws.Range("E9").Value = 0.0001
$ws.Range("F9").Value = 10
$ws.Range("G9").Value = 0.87878787878787801
$ws.Range("H9").Value = 0.86752136752136699
$ws.Range("I9").Value = 0.87311827956989196
$ws.Range("J9").Value = 0.97723143087094499
$ws.Range("D9:J9").Font.Bold = $true

# --- Row 10: Alpha, Beta Attention / Unbalanced ---
$ws.Range("A10").Value = "Alpha, Beta Attention"
$ws.Range("A10").WrapText = $true
$ws.Range("B10").Value = "Alpha, Beta Attention"
$ws.Range("C10").Value = "Unbalanced"
$ws.Range("D10").Value = 450.005
$ws.Range("E10").Value = 0.0001
$ws.Range("F10").Value = 10
$ws.Range("G10").Value = 0.62944162436548201
$ws.Range("H10").Value = 0.48249027237353997
$ws.Range("I10").Value = 0.54625550660792899
$ws.Range("J10").Value = 0.87979205205077504

# --- Row 11: placeholder wrap-text style in column A (kept empty) ---
$ws.Range("A11").WrapText = $true

# --- Selection moves to K14 ---
$ws.Range("K14").Select()

Write-Host "done"
